# Deep Audit complete: Removed all hardcoded defaults
$wb = $excel.ActiveWorkbook

# --- IMPACT_CONFIG sheet ---
$ws1 = $wb.Worksheets.Item("IMPACT_CONFIG")
$ws1.Range("B4").Value = 0    # CO2 TAX RATE ($/Ton)

$ws1.Range("B8").Value = 0    # Solar PV Panels - Unit Cost
$ws1.Range("C8").Value = 0    # Solar PV Panels - CO2 Reduction

$ws1.Range("B9").Value = 0    # Trees Planted - Unit Cost
$ws1.Range("C9").Value = 0    # Trees Planted - CO2 Reduction

$ws1.Range("B10").Value = 0   # Green Electricity - Unit Cost
$ws1.Range("C10").Value = 0   # Green Electricity - CO2 Reduction

$ws1.Range("B11").Value = 0   # CO2 Credits - Unit Cost
$ws1.Range("C11").Value = 0   # CO2 Credits - CO2 Reduction

# --- STRATEGY_SELECTOR sheet ---
$ws2 = $wb.Worksheets.Item("STRATEGY_SELECTOR")
$ws2.Range("B6").Value = 0    # Current CO2 Emissions (Tons/Year)
$ws2.Range("B8").Value = 0    # Energy Consumption (kWh/Year)

$ws2.Range("B13").Value = 0   # Solar PV Panels - Quantity
$ws2.Range("B14").Value = 0   # Trees Planted - Quantity
$ws2.Range("B15").Value = 0   # Green Electricity - Quantity
$ws2.Range("B16").Value = 0   # CO2 Credits - Quantity
